# Update PvsI respirometry rates: recomputed volume/area/rate columns after
# switching the specific-rate area unit from cm2 to m2 (model fitting refit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @{
    2  = @{ T = 0.1450048780487805;  V = 0.0002448603057459146; Z = -0.2487194325456599; AB = -1015.760524303804; AD = -1015.760524303804 }
    3  = @{ T = 0.1492487804878049;  V = 0.0001488973818309612; Z = -0.2716823259261335; AB = -1824.627959103851; AD = -1824.627959103851 }
    4  = @{ T = 0.1469268292682927;  V = 0.0002222807942365138; Z = -0.1793615222778885; AB = -806.9141686035284; AD = -806.9141686035284 }
    5  = @{ T = 0.1418926829268293;  V = 0.0002529432437181515; Z = -0.2709094781434558; AB = -1071.028718384444; AD = -1071.028718384444 }
    6  = @{ T = 0.1446439024390244;  V = 0.0001851607801792304; Z = -0.3107962399977305; AB = -1678.520903275999; AD = -1678.520903275999 }
    7  = @{ T = 0.1429268292682927;  V = 0.0003232296608680373; Z = -0.2284405941259857; AB = -706.7439093074119; AD = -706.7439093074119 }
    9  = @{ T = 0.1450048780487805;  V = 0.0002448603057459146; Z = 0.2909704397192847;  AB = 1188.312000317509;  AD = 1188.312000317509 }
    10 = @{ T = 0.1492487804878049;  V = 0.0001488973818309612; Z = 0.3294622653946687;  AB = 2212.680044090349;  AD = 2212.680044090349 }
    11 = @{ T = 0.1469268292682927;  V = 0.0002222807942365138; Z = 0.1881365243198725;  AB = 846.391272651695;   AD = 846.391272651695 }
    12 = @{ T = 0.1418926829268293;  V = 0.0002529432437181515; Z = 0.3281816441710751;  AB = 1297.451710300513;  AD = 1297.451710300513 }
    13 = @{ T = 0.1446439024390244;  V = 0.0001851607801792304; Z = 0.4119087632202087;  AB = 2224.60049488608;   AD = 2224.60049488608 }
    14 = @{ T = 0.1429268292682927;  V = 0.0003232296608680373; Z = 0.2813515295124803;  AB = 870.4384639606008;  AD = 870.4384639606008 }
}

foreach ($r in $rowUpdates.Keys) {
    $vals = $rowUpdates[$r]
    $ws.Cells.Item($r, 20).Value = $vals.T   # column T = volume
    $ws.Cells.Item($r, 22).Value = $vals.V   # column V = area
    $ws.Cells.Item($r, 26).Value = $vals.Z   # column Z = rate.abs
    $ws.Cells.Item($r, 28).Value = $vals.AB  # column AB = rate.a.spec
    $ws.Cells.Item($r, 30).Value = $vals.AD  # column AD = rate.output
    $ws.Cells.Item($r, 29).Value = "umolO2/min/m2"  # column AC = output.unit
}

# Row 8: area becomes 0, so the per-area specific rate is infinite.
$ws.Cells.Item(8, 20).Value = 0.1544
$ws.Cells.Item(8, 22).Value = 0
$ws.Cells.Item(8, 26).Value = 0.001102594654015669
$ws.Cells.Item(8, 28).Value = "Inf"
$ws.Cells.Item(8, 30).Value = "Inf"
$ws.Cells.Item(8, 29).Value = "umolO2/min/m2"

# Row 15: area is 0 and the absolute rate is 0, so the specific rate is
# undefined (0/0) and the cells are cleared rather than holding a value.
$ws.Cells.Item(15, 20).Value = 0.1544
$ws.Cells.Item(15, 22).Value = 0
$ws.Cells.Item(15, 26).Value = 0
$ws.Cells.Item(15, 28).ClearContents()
$ws.Cells.Item(15, 30).ClearContents()
$ws.Cells.Item(15, 29).Value = "umolO2/min/m2"
